# Adds a new "wx_units" / "metric" key-value row (row 7) to every site sheet
# in the workbook, just above the existing "ET_method" row, shifting the
# "ET_method"/"fao" row down from row 7/8 to row 8/9.
#
# Mirrors the author's commit "Added option to define weather data units."

$wb = $excel.ActiveWorkbook

$siteSheets = @(
    'Agua Fria',
    'Bill Williams 1',
    'Bill Williams 2',
    'Butler',
    'Cienega Creek',
    'Gila Bend',
    'Lower Gila 1',
    'Lower Gila 2',
    'Douglas INA',
    'Duncan Valley',
    'Harquahala INA',
    'Hualapai',
    'Kanab Plateau',
    'Lake Mohave',
    'Little Colorado River 1',
    'Little Colorado River 2',
    'Little Colorado River 3',
    'McMullen Valley',
    'Parker',
    'Phoenix AMA 1',
    'Phoenix AMA 2',
    'Phoenix AMA 3',
    'Pinal AMA 1',
    'Pinal AMA 2',
    'Prescott AMA',
    'Ranegras Plain',
    'Sacramento Valley',
    'Safford',
    'San Simon 1',
    'San Simon 2',
    'San Simon 3',
    'Lower San Pedro',
    'Upper San Pedro',
    'Tonto Creek',
    'Tucson AMA 1',
    'Tucson AMA 2',
    'Verde River 1',
    'Verde River 2',
    'Verde River 3',
    'Virgin River',
    'Willcox 1',
    'Willcox 2',
    'Yuma'
)

foreach ($sheetName in $siteSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Push the "ET_method" row (and everything below it) down one row, and
    # insert a fresh row 7 that inherits formatting from row 6 above it -
    # exactly what Excel's own Insert does.
    $ws.Rows.Item(7).Insert() | Out-Null

    $ws.Range("A7").Value = "wx_units"
    $ws.Range("B7").Value = "metric"

    $ws.Range("A8").Select() | Out-Null
}

$wb.Worksheets.Item("Bill Williams 1").Activate() | Out-Null
